$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 12:05"

# 2. Estados Unidos (row 4) - updated case counts
$ws.Range("B4").Value = 1527951
$ws.Range("C4").Value = 287
$ws.Range("E4").Value = 1090582
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 90980

# 3. Suiza (row 26) - updated case counts
$ws.Range("B26").Value = 30597
$ws.Range("C26").Value = 10
$ws.Range("E26").Value = 1216

# 4. Afganistan (row 59, before re-sort) - updated case counts which push it
#    above Australia/Argelia/Barein/Malasia/Marruecos in the ranking
$ws.Range("B59").Value = 7072
$ws.Range("C59").Value = 408
$ws.Range("D59").Value = 801
$ws.Range("E59").Value = 6098
$ws.Range("G59").Value = 4
$ws.Range("H59").Value = 173

# 5. Hong Kong (row 102) - updated case counts
$ws.Range("D102").Value = 1025
$ws.Range("E102").Value = 27

# 6. Albania (row 107) - updated case counts
$ws.Range("B107").Value = 948
$ws.Range("C107").Value = 2
$ws.Range("D107").Value = 727
$ws.Range("E107").Value = 190

# 7. Re-sort the country table (rows 4-219) by "Casos totales" (column B) descending,
#    matching the sheet's existing sort order, so Afganistan's updated total moves it
#    into its new ranked position.
$sortRange = $ws.Range("A4:H219")
$sortRange.Sort($ws.Range("B4:B219"), 2)
